$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 0.5

# --- Row 5 --- (placeholder result text, drop breakdown columns F:I)
$ws.Range("D5").Value = "----"
$ws.Range("F5:I5").Clear()

# --- Row 6 ---
$ws.Range("D6").Value = "----"
$ws.Range("F6:I6").Clear()

# --- Row 7 ---
$ws.Range("E7").Value = 1

# --- Row 9 ---
$ws.Range("D9").Value = "----"
$ws.Range("E9").Value = 1
$ws.Range("I9").Clear()

# --- Row 10 ---
$ws.Range("D10").Value = "----"
$ws.Range("E10").Value = 1
$ws.Range("I10").Clear()

# --- Row 11 ---
$ws.Range("E11").Value = 0.5

# --- Row 12 ---
$ws.Range("D12").Value = "----"
$ws.Range("E12").Value = 0.5
$ws.Range("I12").Clear()

# --- Row 13 ---
$ws.Range("D13").Value = "----"
$ws.Range("E13").Value = 0.5
$ws.Range("I13").Clear()

# --- Row 14 ---
$ws.Range("D14").Value = "----"
$ws.Range("E14").Value = 0.5
$ws.Range("F14:I14").Clear()

# --- Row 15 ---
$ws.Range("D15").Value = "----"
$ws.Range("E15").Value = 0.5
$ws.Range("F15:I15").Clear()

# --- Row 16 ---
$ws.Range("D16").Value = "----"
$ws.Range("E16").Value = 0.5
$ws.Range("I16").Clear()

# --- Row 17 ---
$ws.Range("D17").Value = "----"
$ws.Range("E17").Value = 0.5
$ws.Range("I17").Clear()

# --- Row 18 ---
$ws.Range("D18").Value = "----"
$ws.Range("E18").Value = 0.5
$ws.Range("F18:I18").Clear()

# --- Row 19 ---
$ws.Range("D19").Value = "----"
$ws.Range("E19").Value = 0.5
$ws.Range("I19").Clear()

# --- Row 20 ---
$ws.Range("D20").Value = "----"
$ws.Range("E20").Value = 0.5
$ws.Range("F20:I20").Clear()

# --- Row 21 ---
$ws.Range("D21").Value = "----"
$ws.Range("E21").Value = 0.5
$ws.Range("F21:I21").Clear()

# --- Row 22 --- (quote-prefixed entry -> gets its own centered/quotePrefix style)
$ws.Range("D22").Value = "'----"
$ws.Range("E22").Value = 0.5
$ws.Range("F22:I22").Clear()

# --- Row 23 ---
$ws.Range("D23").Value = "----"
$ws.Range("E23").Value = 0.5
$ws.Range("F23:I23").Clear()

# --- Row 24 ---
$ws.Range("D24").Value = "----"
$ws.Range("E24").Value = 0.5
$ws.Range("F24:I24").Clear()

# --- Row 25 ---
$ws.Range("D25").Value = "----"
$ws.Range("E25").Value = 0.5
$ws.Range("F25:I25").Clear()

# --- Row 26 ---
$ws.Range("D26").Value = "----"
$ws.Range("E26").Value = 0.5
$ws.Range("F26:I26").Clear()

# --- Row 27 ---
$ws.Range("D27").Value = "----"
$ws.Range("E27").Value = 0.5
$ws.Range("F27:I27").Clear()

# --- Row 28 ---
$ws.Range("D28").Value = "----"
$ws.Range("E28").Value = 0.5
$ws.Range("F28:I28").Clear()

# --- Row 29 ---
$ws.Range("E29").Value = 0.5

# --- Row 30 --- (totals reset to 0 since example results removed)
$ws.Range("E30").Value = 0
$ws.Range("I30").Value = 0

# Restore the selection the author ended up with
$ws.Range("E3:E10").Select()
